$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-05-27T19:07:23"
$ws.Range("W4").Value = 88.56
$ws.Range("X4").Value = 50.79
$ws.Range("Y4").Value = 14.6
$ws.Range("Z4").Value = 6.41
$ws.Range("W6").Value = -1.59
$ws.Range("X6").Value = -0.24
$ws.Range("Y6").Value = 0.01
$ws.Range("Z6").Value = 0.07000000000000001
$ws.Range("X8").Value = 17.07
$ws.Range("W9").Value = 86.94
$ws.Range("X9").Value = 50.96
$ws.Range("Y9").Value = 14.76
$ws.Range("Z9").Value = 6.61
$ws.Range("W11").Value = -3.22
$ws.Range("X11").Value = -0.07000000000000001
$ws.Range("Z11").Value = 0.27
$ws.Range("X13").Value = 17.07
$ws.Range("W14").Value = 88.23999999999999
$ws.Range("X14").Value = 51.4
$ws.Range("Y14").Value = 14.88
$ws.Range("Z14").Value = 6.61
$ws.Range("W15").Value = 1.3
$ws.Range("X15").Value = 0.45
$ws.Range("Y15").Value = 0.12
$ws.Range("W16").Value = -3.22
$ws.Range("X16").Value = -0.07000000000000001
$ws.Range("Z16").Value = 0.27
$ws.Range("X18").Value = 17.07
$ws.Range("W19").Value = 88.65000000000001
$ws.Range("X19").Value = 50.89
$ws.Range("Y19").Value = 14.73
$ws.Range("Z19").Value = 6.5
$ws.Range("W21").Value = -1.51
$ws.Range("X21").Value = -0.14
$ws.Range("Y21").Value = 0.15
$ws.Range("Z21").Value = 0.16
$ws.Range("X23").Value = 17.07
$ws.Range("W24").Value = 88.65000000000001
$ws.Range("X24").Value = 50.89
$ws.Range("Y24").Value = 14.73
$ws.Range("Z24").Value = 6.5
$ws.Range("W26").Value = -1.51
$ws.Range("X26").Value = -0.14
$ws.Range("Y26").Value = 0.15
$ws.Range("Z26").Value = 0.16
$ws.Range("X28").Value = 17.07
$ws.Range("W29").Value = 89
$ws.Range("X29").Value = 51.06
$ws.Range("Y29").Value = 14.91
$ws.Range("Z29").Value = 6.59
$ws.Range("W31").Value = -1.16
$ws.Range("X31").Value = 0.03
$ws.Range("Y31").Value = 0.33
$ws.Range("X33").Value = 17.07
$ws.Range("W34").Value = 86.67
$ws.Range("X34").Value = 124.57
$ws.Range("Y34").Value = 14.96
$ws.Range("Z34").Value = 6.74
$ws.Range("W35").Value = 1.3
$ws.Range("X35").Value = 0.45
$ws.Range("Y35").Value = 0.12
$ws.Range("W36").Value = -4.78
$ws.Range("X36").Value = -0.07000000000000001
$ws.Range("Y36").Value = 0.25
$ws.Range("Z36").Value = 0.4
$ws.Range("X37").Value = 73.17
$ws.Range("X38").Value = 17.07
$ws.Range("W39").Value = 88.56
$ws.Range("X39").Value = 50.79
$ws.Range("Y39").Value = 14.6
$ws.Range("Z39").Value = 6.41
$ws.Range("W41").Value = -1.59
$ws.Range("X41").Value = -0.24
$ws.Range("Y41").Value = 0.01
$ws.Range("Z41").Value = 0.07000000000000001
$ws.Range("X43").Value = 17.07
$ws.Range("W44").Value = 91.25
$ws.Range("X44").Value = 51.51
$ws.Range("Y44").Value = 14.6
$ws.Range("Z44").Value = 6.29
$ws.Range("W46").Value = 1.1
$ws.Range("X46").Value = 0.48
$ws.Range("Z46").Value = -0.05
$ws.Range("X48").Value = 17.07
$ws.Range("W49").Value = 81.22
$ws.Range("X49").Value = 47.74
$ws.Range("Y49").Value = 12.97
$ws.Range("Z49").Value = 5.51
$ws.Range("W51").Value = -8.93
$ws.Range("X51").Value = -3.28
$ws.Range("Y51").Value = -1.61
$ws.Range("Z51").Value = -0.84
$ws.Range("X53").Value = 17.07
$ws.Range("W54").Value = 79.64
$ws.Range("X54").Value = 47.12
$ws.Range("Y54").Value = 12.84
$ws.Range("Z54").Value = 5.59
$ws.Range("W56").Value = -10.51
$ws.Range("X56").Value = -3.91
$ws.Range("Y56").Value = -1.75
$ws.Range("Z56").Value = -0.75
$ws.Range("X58").Value = 17.07
$ws.Range("W59").Value = 93.52
$ws.Range("X59").Value = 52.33
$ws.Range("Y59").Value = 14.85
$ws.Range("Z59").Value = 6.36
$ws.Range("W61").Value = 3.37
$ws.Range("X61").Value = 1.3
$ws.Range("Z61").Value = 0.02
$ws.Range("X63").Value = 17.07
$ws.Range("W64").Value = 95.09999999999999
$ws.Range("X64").Value = 52.81
$ws.Range("Y64").Value = 15.03
$ws.Range("Z64").Value = 6.43
$ws.Range("W66").Value = 4.95
$ws.Range("X66").Value = 1.79
$ws.Range("Y66").Value = 0.45
$ws.Range("Z66").Value = 0.08
$ws.Range("X68").Value = 17.07
$ws.Range("W69").Value = 95.70999999999999
$ws.Range("X69").Value = 53.27
$ws.Range("Y69").Value = 15.08
$ws.Range("Z69").Value = 6.39
$ws.Range("W71").Value = 5.55
$ws.Range("X71").Value = 2.24
$ws.Range("Y71").Value = 0.5
$ws.Range("Z71").Value = 0.05
$ws.Range("X73").Value = 17.07
$ws.Range("W74").Value = 91.81
$ws.Range("X74").Value = 51.61
$ws.Range("Y74").Value = 14.55
$ws.Range("Z74").Value = 6.23
$ws.Range("W76").Value = 1.65
$ws.Range("X76").Value = 0.59
$ws.Range("Z76").Value = -0.11
$ws.Range("X78").Value = 17.07
$ws.Range("W79").Value = 90.16
$ws.Range("X79").Value = 51.03
$ws.Range("Y79").Value = 14.58
$ws.Range("Z79").Value = 6.34
$ws.Range("X83").Value = 17.07
$ws.Range("W84").Value = 79.70999999999999
$ws.Range("X84").Value = 47.12
$ws.Range("Y84").Value = 12.85
$ws.Range("Z84").Value = 5.58
$ws.Range("W86").Value = -10.44
$ws.Range("X86").Value = -3.91
$ws.Range("Y86").Value = -1.73
$ws.Range("Z86").Value = -0.76
$ws.Range("X88").Value = 17.07
$ws.Range("W89").Value = 89
$ws.Range("X89").Value = 51.06
$ws.Range("Y89").Value = 14.91
$ws.Range("Z89").Value = 6.59
$ws.Range("W91").Value = -1.16
$ws.Range("X91").Value = 0.03
$ws.Range("Y91").Value = 0.33
$ws.Range("X93").Value = 17.07
